$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Date Entered?" (G) and "Time Entered?" (H) columns entirely.
$ws.Range("G1:H3").Delete()

# New rows of check_availability results (rows 4-10).
$rows = @(
    @("2024-09-09 15:28:31", "check_availability", "https://www.opentable.com/r/bar-spero-washington/", "No availability for the selected date.", "2024-09-09", "15:28:31"),
    @("2024-09-09 15:34:29", "check_availability", "https://www.opentable.com/r/bar-spero-washington/", "No availability for the selected date.", "2024-09-09", "15:34:29"),
    @("2024-09-09 15:38:11", "check_availability", "https://www.opentable.com/r/bar-spero-washington/", "No availability for the selected date.", "2024-09-09", "15:38:11"),
    @("2024-09-09 15:40:38", "check_availability", "https://www.opentable.com/r/bar-spero-washington/", "No availability for the selected date.", "2024-09-09", "15:40:38"),
    @("2024-09-09 15:43:17", "check_availability", "https://www.opentable.com/r/bar-spero-washington/", "No availability for the selected date.", "2024-09-09", "15:43:17"),
    @("2024-09-09 15:44:02", "check_availability", "https://www.opentable.com/r/bar-spero-washington/", "No availability for the selected date.", "2024-09-09", "15:44:02"),
    @("2024-09-09 15:45:01", "check_availability", "https://www.opentable.com/r/bar-spero-washington/", "No availability for the selected date.", "2024-09-09", "15:45:01")
)

$r = 4
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]

    # Columns E (date) and F (time) would otherwise be auto-converted to
    # date/time serial numbers by Excel, so force them to be entered as
    # text, then strip the formatting Excel applied so the cell keeps the
    # default (no explicit) style, matching plain text cells elsewhere.
    $ws.Cells.Item($r, 5).NumberFormat = "@"
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 5).ClearFormats()

    $ws.Cells.Item($r, 6).NumberFormat = "@"
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 6).ClearFormats()

    $r = $r + 1
}
